$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and restore text-number formatting
# so strings like "1.00" are not auto-converted to numeric values by Excel).

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '90.163.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range('E2').Value = '  +4.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '3.231.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '219.61'
$ws.Range("E5").NumberFormat = "@"
$ws.Range('E5').Value = '  +6.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '634.07'
$ws.Range("E6").NumberFormat = "@"
$ws.Range('E6').Value = '  +4.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.399'
$ws.Range("E7").NumberFormat = "@"
$ws.Range('E7').Value = '  +8.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.703'
$ws.Range("E8").NumberFormat = "@"
$ws.Range('E8').Value = '  +7.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range("E9").NumberFormat = "@"
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '3.230.00'
$ws.Range("E10").NumberFormat = "@"
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.583'
$ws.Range("E11").NumberFormat = "@"
$ws.Range('E11').Value = '  +9.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '0.182'
$ws.Range("E12").NumberFormat = "@"
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '0.0000266'
$ws.Range("E13").NumberFormat = "@"
$ws.Range('E13').Value = '  +10.22%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range('B14').Value = 'Avalanche'
$ws.Range("C14").NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '34.27'
$ws.Range("E14").NumberFormat = "@"
$ws.Range('E14').Value = '  +7.36%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range('B15').Value = 'Toncoin'
$ws.Range("C15").NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '5.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range('E15').Value = '  +4.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '3.831.05'
$ws.Range("E16").NumberFormat = "@"
$ws.Range('E16').Value = '  +2.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '89.951.03'
$ws.Range("E17").NumberFormat = "@"
$ws.Range('E17').Value = '  +4.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '3.228.25'
$ws.Range("E18").NumberFormat = "@"
$ws.Range('E18').Value = '  +3.29%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range("C19").NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '3.48'
$ws.Range("E19").NumberFormat = "@"
$ws.Range('E19').Value = '  +17.48%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range('B20').Value = 'PEPE'
$ws.Range("C20").NumberFormat = "@"
$ws.Range('C20').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '0.0000233'
$ws.Range("E20").NumberFormat = "@"
$ws.Range('E20').Value = '  +80.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '13.79'
$ws.Range("E21").NumberFormat = "@"
$ws.Range('E21').Value = '  +3.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '445.15'
$ws.Range("E22").NumberFormat = "@"
$ws.Range('E22').Value = '  +8.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '8.77'
$ws.Range("E23").NumberFormat = "@"
$ws.Range('E23').Value = '  +4.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '5.18'
$ws.Range("E24").NumberFormat = "@"
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range('B25').Value = 'Aptos'
$ws.Range("C25").NumberFormat = "@"
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '12.20'
$ws.Range("E25").NumberFormat = "@"
$ws.Range('E25').Value = '  +4.18%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range("C26").NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '5.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range('E26').Value = '  +3.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '83.74'
$ws.Range("E27").NumberFormat = "@"
$ws.Range('E27').Value = '  +14.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '3.417.40'
$ws.Range("E28").NumberFormat = "@"
$ws.Range('E28').Value = '  +2.79%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '0.163'
$ws.Range("E30").NumberFormat = "@"
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '4.17'
$ws.Range("E32").NumberFormat = "@"
$ws.Range('E32').Value = '  +39.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '8.68'
$ws.Range("E33").NumberFormat = "@"
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '553.66'
$ws.Range("E34").NumberFormat = "@"
$ws.Range('E34').Value = '  +3.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '7.13'
$ws.Range("E35").NumberFormat = "@"
$ws.Range('E35').Value = '  +8.57%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range("C36").NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '1.35'
$ws.Range("E36").NumberFormat = "@"
$ws.Range('E36').Value = '  +6.42%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range("C37").NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '1.93'
$ws.Range("E37").NumberFormat = "@"
$ws.Range('E37').Value = '  +4.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '22.69'
$ws.Range("E38").NumberFormat = "@"
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range('B39').Value = 'Kaspa'
$ws.Range("C39").NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.132'
$ws.Range("E39").NumberFormat = "@"
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range("C40").NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '22.42'
$ws.Range("E40").NumberFormat = "@"
$ws.Range('E40').Value = '  +3.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range("E41").NumberFormat = "@"
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '1.97'
$ws.Range("E42").NumberFormat = "@"
$ws.Range('E42').Value = '  +4.44%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '0.381'
$ws.Range("E43").NumberFormat = "@"
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range('B44').Value = 'USDe'
$ws.Range("C44").NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range("E44").NumberFormat = "@"
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '147.12'
$ws.Range("E45").NumberFormat = "@"
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '175.94'
$ws.Range("E46").NumberFormat = "@"
$ws.Range('E46').Value = '  +2.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '44.19'
$ws.Range("E47").NumberFormat = "@"
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.768'
$ws.Range("E48").NumberFormat = "@"
$ws.Range('E48').Value = '  +12.29%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range('B49').Value = 'Stellar'
$ws.Range("C49").NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '0.126'
$ws.Range("E49").NumberFormat = "@"
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range("C50").NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '1.27'
$ws.Range("E50").NumberFormat = "@"
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '0.632'
$ws.Range("E51").NumberFormat = "@"
$ws.Range('E51').Value = '  +7.96%  '
